# Apply "updated figure and data" changes to the NOx sheet:
#  - Add "Figure 3" label (col H) and source-URL hyperlink (col G) to rows 3-5
#  - Add a brand new row 6 ("Net NOx Emissions") with the same structure
#  - Add matching hyperlinks for every new G-cell

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NOx")

$sourceUrl = "https://doi.org/10.1016/j.atmosenv.2020.117834"

# --- New row of data (row 6): Net NOx Emissions ---------------------------
$ws.Range("A6").Value = 17.5
$ws.Range("B6").Value = 0.6
$ws.Range("C6").Value = 29
$ws.Range("D6").Value = "Lee et al."
$ws.Range("E6").Value = 2018
$ws.Range("F6").Value = "Net NOx Emissions"

# --- Fill in the "Figure 3" reference column (H) for rows 3-6 -------------
$ws.Range("H3").Value = "Figure 3"
$ws.Range("H4").Value = "Figure 3"
$ws.Range("H5").Value = "Figure 3"
$ws.Range("H6").Value = "Figure 3"

# --- Add source-URL hyperlinks to column G for rows 3-6 --------------------
# (Hyperlinks.Add sets the cell's text/value; re-applying the Hyperlink
#  style afterwards keeps it identical to the one already used in G2.)
foreach ($r in 3..6) {
    $cell = $ws.Range("G$r")
    [void]$ws.Hyperlinks.Add($cell, $sourceUrl)
    $cell.Style = "Hyperlink"
}
